# Auto-generated Excel COM-interop script
# Applies numeric corrections to H:N columns across multiple sheets
# as described by the commit diff for Mandragora_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 982.65
$ws.Range("I88").Value = 564.125
$ws.Range("J88").Value = 1261.6666
$ws.Range("K88").Value = 564.125
$ws.Range("L88").Value = 1261.6666
$ws.Range("M88").Value = -158.125
$ws.Range("N88").Value = -2073.6666

$ws.Range("H91").Value = 982.65
$ws.Range("I91").Value = 564.125
$ws.Range("J91").Value = 1261.6666
$ws.Range("K91").Value = 564.125
$ws.Range("L91").Value = 1261.6666
$ws.Range("M91").Value = 839.875
$ws.Range("N91").Value = -4069.6666

$ws.Range("H94").Value = 3801.2
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 5002
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 5002
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -5904

$ws.Range("H97").Value = 3343.25
$ws.Range("J97").Value = 3343.25
$ws.Range("L97").Value = 10029.75
$ws.Range("N97").Value = -11021.75

$ws.Range("H98").Value = 1532.8667
$ws.Range("I98").Value = 1470.1
$ws.Range("K98").Value = 1470.1
$ws.Range("M98").Value = 27.90000000000009

$ws.Range("H100").Value = 2678
$ws.Range("I100").Value = 1975
$ws.Range("J100").Value = 5490
$ws.Range("K100").Value = 1975
$ws.Range("L100").Value = 5490
$ws.Range("M100").Value = -1434
$ws.Range("N100").Value = -6572

$ws.Range("H122").Value = 1532.8667
$ws.Range("I122").Value = 1470.1
$ws.Range("K122").Value = 4410.299999999999
$ws.Range("M122").Value = -1960.299999999999

$ws.Range("H132").Value = 5315.697
$ws.Range("I132").Value = 4641.104
$ws.Range("J132").Value = 7114.6113
$ws.Range("K132").Value = 13923.312
$ws.Range("L132").Value = 21343.8339
$ws.Range("M132").Value = -11393.312
$ws.Range("N132").Value = -26403.8339

$ws.Range("H135").Value = 999.05
$ws.Range("I135").Value = 681.2353000000001
$ws.Range("J135").Value = 2800
$ws.Range("K135").Value = 6131.117700000001
$ws.Range("L135").Value = 25200
$ws.Range("M135").Value = -3596.117700000001
$ws.Range("N135").Value = -30270

$ws.Range("H137").Value = 1882.25
$ws.Range("I137").Value = 3689.2222
$ws.Range("J137").Value = 1026.3158
$ws.Range("K137").Value = 11067.6666
$ws.Range("L137").Value = 3078.9474
$ws.Range("M137").Value = -8517.6666
$ws.Range("N137").Value = -8178.9474

$ws.Range("H138").Value = 1536.2826
$ws.Range("I138").Value = 1162.3478
$ws.Range("J138").Value = 1910.2174
$ws.Range("K138").Value = 3487.0434
$ws.Range("L138").Value = 5730.6522
$ws.Range("M138").Value = 1652.9566
$ws.Range("N138").Value = -16010.6522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9843.237999999999
$ws.Range("I32").Value = 8468.258
$ws.Range("J32").Value = 13718.182
$ws.Range("K32").Value = 8468.258
$ws.Range("L32").Value = 13718.182
$ws.Range("M32").Value = -8181.258
$ws.Range("N32").Value = -14292.182

$ws.Range("H61").Value = 1324.8125
$ws.Range("I61").Value = 1075.5652
$ws.Range("J61").Value = 1961.7778
$ws.Range("K61").Value = 1075.5652
$ws.Range("L61").Value = 1961.7778
$ws.Range("M61").Value = -863.5652
$ws.Range("N61").Value = -2385.7778

$ws.Range("H132").Value = 6297.6763
$ws.Range("I132").Value = 4599.6665
$ws.Range("J132").Value = 7223.864
$ws.Range("K132").Value = 13798.9995
$ws.Range("L132").Value = 21671.592
$ws.Range("M132").Value = -11268.9995
$ws.Range("N132").Value = -26731.592

$ws.Range("H136").Value = 1324.8125
$ws.Range("I136").Value = 1075.5652
$ws.Range("J136").Value = 1961.7778
$ws.Range("K136").Value = 3226.6956
$ws.Range("L136").Value = 5885.3334
$ws.Range("M136").Value = -676.6956
$ws.Range("N136").Value = -10985.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4685.137
$ws.Range("I134").Value = 2197.25
$ws.Range("J134").Value = 6896.593
$ws.Range("K134").Value = 6591.75
$ws.Range("L134").Value = 20689.779
$ws.Range("M134").Value = -4056.75
$ws.Range("N134").Value = -25759.779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2288.3635
$ws.Range("I122").Value = 1837.5
$ws.Range("J122").Value = 2546
$ws.Range("K122").Value = 5512.5
$ws.Range("L122").Value = 7638
$ws.Range("M122").Value = -3062.5
$ws.Range("N122").Value = -12538

$ws.Range("H132").Value = 2627.1843
$ws.Range("I132").Value = 1948.579
$ws.Range("J132").Value = 3305.7896
$ws.Range("K132").Value = 5845.737
$ws.Range("L132").Value = 9917.3688
$ws.Range("M132").Value = -3315.737
$ws.Range("N132").Value = -14977.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1016.4074
$ws.Range("I98").Value = 662.5
$ws.Range("J98").Value = 1224.5883
$ws.Range("K98").Value = 1987.5
$ws.Range("L98").Value = 3673.7649
$ws.Range("M98").Value = -489.5
$ws.Range("N98").Value = -6669.7649

$ws.Range("H122").Value = 2699.0938
$ws.Range("I122").Value = 453.14285
$ws.Range("J122").Value = 3327.96
$ws.Range("K122").Value = 4078.28565
$ws.Range("L122").Value = 29951.64
$ws.Range("M122").Value = -1628.28565
$ws.Range("N122").Value = -34851.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2101.6667
$ws.Range("I97").Value = 1870
$ws.Range("J97").Value = 2333.3333
$ws.Range("K97").Value = 1870
$ws.Range("L97").Value = 2333.3333
$ws.Range("M97").Value = -1374
$ws.Range("N97").Value = -3325.3333

$ws.Range("H107").Value = 691.25806
$ws.Range("I107").Value = 705.2222
$ws.Range("J107").Value = 671.9231
$ws.Range("K107").Value = 705.2222
$ws.Range("L107").Value = 671.9231
$ws.Range("M107").Value = 1214.7778
$ws.Range("N107").Value = -4511.9231

$ws.Range("H122").Value = 3723.9
$ws.Range("I122").Value = 2981
$ws.Range("J122").Value = 4838.25
$ws.Range("K122").Value = 8943
$ws.Range("L122").Value = 14514.75
$ws.Range("M122").Value = -6493
$ws.Range("N122").Value = -19414.75

$ws.Range("H132").Value = 2328
$ws.Range("I132").Value = 1883.9445
$ws.Range("J132").Value = 2675.5217
$ws.Range("K132").Value = 5651.833500000001
$ws.Range("L132").Value = 8026.5651
$ws.Range("M132").Value = -3121.833500000001
$ws.Range("N132").Value = -13086.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3533
$ws.Range("I122").Value = 686.6667
$ws.Range("K122").Value = 2060.0001
$ws.Range("M122").Value = 389.9998999999998

$ws.Range("H126").Value = 2250.1667
$ws.Range("I126").Value = 2465.4783
$ws.Range("K126").Value = 7396.4349
$ws.Range("M126").Value = -4926.4349
